$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric values in D2:D9 with text values "123-1" .. "130-8"
$values = @("123-1", "124-2", "125-3", "126-4", "127-5", "128-6", "129-7", "130-8")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Update the active selection to G6
$ws.Range("G6").Select()
